# Commit Board and Card
#
# Original workbook:
#   Tab 1 = "CreateACard" (sheetId=2)  -> big TESTCASE table (A1:K4)
#   Tab 2 = "Create"      (sheetId=1)  -> small 2-col table (A1:B4)
#
# Target workbook:
#   Tab 1 = "Board" (sheetId=1) -> brand new small TESTCASE table (A1:I6)
#   Tab 2 = "Card"  (sheetId=2) -> same big TESTCASE table, values updated
#
# Because sheetId follows sheet *identity* (not tab position) in this
# engine, we move the "Create" sheet (sheetId=1) to tab 1 and turn it into
# "Board", and let "CreateACard" (sheetId=2) settle into tab 2 as "Card".

$wb = $excel.ActiveWorkbook

$sheetCreateACard = $wb.Worksheets.Item(1)   # sheetId=2, big table
$sheetCreate      = $wb.Worksheets.Item(2)   # sheetId=1, small table

# --- Reorder: put "Create" (sheetId=1) first -------------------------------
$sheetCreate.Move($sheetCreateACard)

$board = $wb.Worksheets.Item(1)   # was "Create", sheetId=1
$card  = $wb.Worksheets.Item(2)   # was "CreateACard", sheetId=2

$board.Name = "Board"
$card.Name  = "Card"

# =============================================================================
# Shared string literals (kept as variables to avoid transcription mistakes)
# =============================================================================
$s_name             = "name"
$s_no               = "No."
$s_testcase         = "TESTCASE"
$s_method           = "METHOD"
$s_post             = "POST"
$s_put              = "PUT"
$s_status_code      = "STATUS_CODE"
$s_status_message   = "STATUS_MESSAGE"
$s_desc             = "desc"
$s_start            = "start"
$s_minitype         = "miniType"
$s_high             = "High"
$s_create_card_ok   = "Create a Card successfully"
$s_bug              = "This is a bug"
$s_create_board_ok  = "Create A Board Successfully"
$s_update_board_ok  = "Update A Board Successfully"
$s_new_chapter      = "The new Chapter"
$s_upd_new_chapter  = "[Updated] The new Chapter"
$s_upd_new_bug      = "[Updated] New Bug"
$s_date1            = "2023-03-29"
$s_date2            = "2023-03-30"
$s_lorem_trail      = "Lorem Ipsum is simply dummy text of the printing and typesetting industry. Lorem Ipsum has been the industry's standard dummy text ever since the 1500s, when an unknown printer took a galley of type and scrambled it to make a type specimen book. "
$s_lorem_notrail    = "Lorem Ipsum is simply dummy text of the printing and typesetting industry. Lorem Ipsum has been the industry's standard dummy text ever since the 1500s, when an unknown printer took a galley of type and scrambled it to make a type specimen book"
$s_upd_lorem_notrail= "[Updated] Lorem Ipsum has been the industry's standard dummy text ever since the 1500s, when an unknown printer took a galley of type and scrambled it to make a type specimen book"
$s_upd_lorem_trail  = "[Updated] Lorem Ipsum has been the industry's standard dummy text ever since the 1500s, when an unknown printer took a galley of type and scrambled it to make a type specimen book. "

# =============================================================================
# BOARD sheet (tab 1) - wipe the old "Create" content and rebuild from scratch
# =============================================================================

# Wipe first (Clear() also drops the clipboard, so do this before Copy()).
$board.Cells.Clear()

# Reuse the existing "header" style (s=1) from the Card sheet's A1 instead of
# re-building fill/alignment by hand (which would create brand-new style
# entries in styles.xml).
$card.Range("A1").Copy()
$board.Range("A1:I1").PasteSpecial(-4122)   # xlPasteFormats

$board.Cells.Item(1,1).Value = $s_no
$board.Cells.Item(1,2).Value = $s_testcase
$board.Cells.Item(1,3).Value = $s_method
$board.Cells.Item(1,4).Value = $s_name
$board.Cells.Item(1,5).Value = $s_desc
$board.Cells.Item(1,6).Value = $s_status_code
$board.Cells.Item(1,7).Value = $s_status_message
$board.Cells.Item(1,8).Value = $s_no
$board.Cells.Item(1,9).Value = $s_no

$board.Cells.Item(2,1).Value = 1
$board.Cells.Item(2,2).Value = $s_create_board_ok
$board.Cells.Item(2,3).Value = $s_post
$board.Cells.Item(2,4).Value = $s_new_chapter
$board.Cells.Item(2,5).Value = $s_lorem_notrail
$board.Cells.Item(2,6).Value = 200

$board.Cells.Item(3,1).Value = 2
$board.Cells.Item(3,2).Value = $s_update_board_ok
$board.Cells.Item(3,3).Value = $s_put
$board.Cells.Item(3,4).Value = $s_upd_new_chapter
$board.Cells.Item(3,5).Value = $s_upd_lorem_notrail
$board.Cells.Item(3,6).Value = 200

$board.Cells.Item(4,1).Value = 3
$board.Cells.Item(5,1).Value = 4
$board.Cells.Item(6,1).Value = 5

$board.Columns.Item(2).ColumnWidth = 24.582589285714285
$board.Columns.Item(4).ColumnWidth = 14.984933035714286
$board.Columns.Item(5).ColumnWidth = 18.883370535714285
$board.Columns.Item(6).ColumnWidth = 17.383370535714285
$board.Columns.Item(7).ColumnWidth = 21.484933035714285

# =============================================================================
# CARD sheet (tab 2) - keep the existing layout/styles, edit values in place
# =============================================================================

# Row 2, col F: date serial -> literal text date (drives numFmt 165 -> 49)
$card.Range("F2").NumberFormat = "@"
$card.Range("F2").Value = $s_date1

# Row 3: fill in the previously-empty cells
$card.Range("B3").Value = $s_create_card_ok
$card.Range("C3").Value = $s_put
$card.Range("D3").Value = $s_upd_new_bug
$card.Range("E3").Value = $s_upd_lorem_trail
$card.Range("F3").NumberFormat = "@"
$card.Range("F3").Value = $s_date2

# New column widths (D and E), inserted between existing B and F widths
$card.Columns.Item(4).ColumnWidth = 15.785714285714286
$card.Columns.Item(5).ColumnWidth = 17.785714285714285

# =============================================================================
# Selections / active tab - activating each sheet + selecting sets both the
# per-sheet <selection> and the workbook's activeTab/tabSelected bookkeeping.
# =============================================================================
$board.Activate()
$board.Range("D16").Select()

$card.Activate()
$card.Range("G9").Select()
